$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: task changed from user-creation-via-website to "possibility to create a user"
$ws.Range("B18").Value = "Möjlighet att skapa användare"
$ws.Range("C18").Value = "Front end"
$ws.Range("E18").Style = "Bra"
$ws.Range("F18").Value = "S2"
$ws.Range("G18").ClearContents()

# Update the active selection to B15
$ws.Range("B15").Select() | Out-Null
